$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.066.87"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.567.81"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "'208.47"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  +0.95%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("D12").Value = "1.566.23"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "'3.78"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "27.060.82"
$ws.Range("E15").Value = "  +0.55%  "
$ws.Range("D16").Value = "'61.93"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("D19").Value = "'215.72"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "'4.16"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").Value = "'9.21"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "'6.63"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").Value = "'3.24"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").Value = "1.423.54"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("E34").Value = "  +12.80%  "
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "'0.533"
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").Value = "'0.813"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  +2.55%  "
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "'64.74"
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "1.703.83"
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "0.0₆0102"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  +0.32%  "
$ws.Range("E51").Value = "  +0.39%  "
